# Added monthly trends to RSTs
# The underlying data rows for the NCEP (rows 4/5), ERA (rows 9/10) and
# ERA 2.5 (rows 14/15) blocks had their "West" and "Central" detail rows
# swapped. Column A (the row label) stays put; columns B:T need to be
# exchanged between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$rowPairs = @(@(4,5), @(9,10), @(14,15))

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v2 -eq $null) {
            $cell1.ClearContents()
        } else {
            $cell1.Value2 = $v2
        }

        if ($v1 -eq $null) {
            $cell2.ClearContents()
        } else {
            $cell2.Value2 = $v1
        }
    }
}
